$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and column E (Volume) cells are stored as literal text
# in the source sheet (prices use "." as a thousands separator in several
# rows, e.g. "26.719.32", so the whole column is text, not numeric).
# For values that Excel would otherwise auto-detect as a number, force the
# cell to Text before writing, then drop back to the default (Normal) style
# so no stray number-format is left behind on the cell.

$ws.Cells.Item(2, 4).Value = '26.719.32'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '
$ws.Cells.Item(3, 4).Value = '1.600.78'
$ws.Cells.Item(3, 5).Value = '  +0.24%  '
$ws.Cells.Item(4, 5).Value = '  +0.24%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '211.30'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.07%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.512'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.46%  '
$ws.Cells.Item(7, 5).Value = '  +0.21%  '
$ws.Cells.Item(8, 5).Value = '  +0.28%  '
$ws.Cells.Item(9, 5).Value = '  +0.40%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.67'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.07%  '
$ws.Cells.Item(11, 5).Value = '  +0.71%  '
$ws.Cells.Item(12, 4).Value = '1.826.17'
$ws.Cells.Item(12, 5).Value = '  +0.30%  '
$ws.Cells.Item(13, 4).Value = '1.602.84'
$ws.Cells.Item(13, 5).Value = '  +0.65%  '
$ws.Cells.Item(14, 5).Value = '  +0.45%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.523'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.26%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '65.20'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.30%  '
$ws.Cells.Item(17, 4).Value = '26.696.26'
$ws.Cells.Item(17, 5).Value = '  +0.31%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0743'
$ws.Cells.Item(18, 5).Value = '  +0.85%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '210.97'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.20'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.41%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.00'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.19%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.30'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.91%  '
$ws.Cells.Item(23, 5).Value = '  +0.15%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.97'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.89%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '143.65'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.08%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.16%  '
$ws.Cells.Item(27, 5).Value = '  -0.31%  '
$ws.Cells.Item(28, 5).Value = '  -0.74%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '15.39'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.94%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0513'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.11%  '
$ws.Cells.Item(31, 5).Value = '  -0.27%  '
$ws.Cells.Item(32, 5).Value = '  +1.25%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.97'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.55%  '
$ws.Cells.Item(34, 4).Value = '1.296.12'
$ws.Cells.Item(34, 5).Value = '  +1.60%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.48'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.87%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.607'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -2.00%  '
$ws.Cells.Item(37, 5).Value = '  +0.98%  '
$ws.Cells.Item(38, 5).Value = '  +20.91%  '
$ws.Cells.Item(39, 5).Value = '  -0.37%  '
$ws.Cells.Item(40, 5).Value = '  -1.84%  '
$ws.Cells.Item(41, 5).Value = '  -1.19%  '
$ws.Cells.Item(42, 5).Value = '  -0.27%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.784'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.16%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '63.23'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -1.24%  '
$ws.Cells.Item(45, 4).Value = '1.737.87'
$ws.Cells.Item(45, 5).Value = '  +0.29%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '91.05'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.10%  '
$ws.Cells.Item(47, 5).Value = '  -2.63%  '
$ws.Cells.Item(48, 5).Value = '  -0.40%  '
$ws.Cells.Item(49, 5).Value = '  +1.84%  '
$ws.Cells.Item(50, 5).Value = '  +0.05%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.44'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.34%  '
